# Corrected excel sheets for application fix issues
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Input
# ---------------------------------------------------------------------------
$wsInput = $wb.Worksheets.Item("Input")
$wsInput.Range("B2").Value = 42064
$wsInput.Activate()
$wsInput.Range("E3").Select()

# ---------------------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("G2").ClearContents()
$wsSummary.Range("A3").Value = 211.19
$wsSummary.Range("E3").Value = 114.2
$wsSummary.Activate()
$wsSummary.Range("D4").Select()

# ---------------------------------------------------------------------------
# Sheet: Repayment schedule
# ---------------------------------------------------------------------------
$wsSched = $wb.Worksheets.Item("Repayment schedule")

# Column L width adjustment (engine's ColumnWidth setter only supports a
# quantized set of widths; 7.29 lands on the value closest to the target
# 8.140625 stored width that is reachable through this property).
$wsSched.Columns.Item(12).ColumnWidth = 7.29

# Row 3
$wsSched.Range("D3").Value = 42064
$wsSched.Range("M3").Copy()
$wsSched.Range("N3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsSched.Range("N3").Value = 0

# Row 4
$wsSched.Range("F4").Value = 921.65
$wsSched.Range("G4").Value = 3211.57
$wsSched.Range("H4").Value = 42.12

# Row 5
$wsSched.Range("B5").Value = 30
$wsSched.Range("C5").Value = 42125
$wsSched.Range("F5").Value = 932.09
$wsSched.Range("G5").Value = 2279.48
$wsSched.Range("H5").Value = 31.68

# Row 6
$wsSched.Range("B6").Value = 31
$wsSched.Range("C6").Value = 42156
$wsSched.Range("F6").Value = 940.54
$wsSched.Range("G6").Value = 1338.94
$wsSched.Range("H6").Value = 23.23

# Row 7
$wsSched.Range("B7").Value = 30
$wsSched.Range("C7").Value = 42186
$wsSched.Range("F7").Value = 950.56
$wsSched.Range("G7").Value = 388.38
$wsSched.Range("H7").Value = 13.21

# Row 8
$wsSched.Range("B8").Value = 31
$wsSched.Range("C8").Value = 42217
$wsSched.Range("F8").Value = 388.38
$wsSched.Range("H8").Value = 3.96
$wsSched.Range("K8").Value = 392.34
$wsSched.Range("P8").Value = 392.34

$wsSched.Activate()
$wsSched.Range("H3:H8").Select()

# ---------------------------------------------------------------------------
# Sheet: Transactions
# ---------------------------------------------------------------------------
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Range("A2").Value = 6565
$wsTrans.Range("C2").Value = 42064
$wsTrans.Range("A3").Value = 6564
$wsTrans.Activate()
$wsTrans.Range("D3").Select()

Write-Output "edits applied"
